# Auto-generated: apply updated market price / profit values per diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 652.4286
$ws.Range("I41").Value = 553.86664
$ws.Range("J41").Value = 898.8333
$ws.Range("K41").Value = 553.86664
$ws.Range("L41").Value = 898.8333
$ws.Range("M41").Value = -113.86664
$ws.Range("N41").Value = -1778.8333
$ws.Range("H43").Value = 1487.375
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 1649.8334
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 1649.8334
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -1787.8334
$ws.Range("H55").Value = 156.88889
$ws.Range("I55").Value = 112.4
$ws.Range("K55").Value = 112.4
$ws.Range("M55").Value = 101.6
$ws.Range("H68").Value = 33294.5
$ws.Range("J68").Value = 33294.5
$ws.Range("L68").Value = 33294.5
$ws.Range("N68").Value = -34792.5
$ws.Range("H71").Value = 33294.5
$ws.Range("J71").Value = 33294.5
$ws.Range("L71").Value = 99883.5
$ws.Range("N71").Value = -107371.5
$ws.Range("H76").Value = 3624.2292
$ws.Range("I76").Value = 3479.389
$ws.Range("K76").Value = 3479.389
$ws.Range("M76").Value = -3164.389
$ws.Range("H79").Value = 3624.2292
$ws.Range("I79").Value = 3479.389
$ws.Range("K79").Value = 3479.389
$ws.Range("M79").Value = -2387.389
$ws.Range("H133").Value = 59320.527
$ws.Range("J133").Value = 59320.527
$ws.Range("L133").Value = 59320.527
$ws.Range("N133").Value = -69440.527
$ws.Range("H137").Value = 3034.8545
$ws.Range("I137").Value = 1751.4286
$ws.Range("J137").Value = 3827.5588
$ws.Range("K137").Value = 5254.2858
$ws.Range("L137").Value = 11482.6764
$ws.Range("M137").Value = -2704.2858
$ws.Range("N137").Value = -16582.6764
$ws.Range("H138").Value = 1216932.4
$ws.Range("I138").Value = 3114.1333
$ws.Range("J138").Value = 1520386.9
$ws.Range("K138").Value = 9342.3999
$ws.Range("L138").Value = 4561160.699999999
$ws.Range("M138").Value = -4202.3999
$ws.Range("N138").Value = -4571440.699999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13900.571
$ws.Range("I32").Value = 11654.372
$ws.Range("J32").Value = 29998.334
$ws.Range("K32").Value = 11654.372
$ws.Range("L32").Value = 29998.334
$ws.Range("M32").Value = -11367.372
$ws.Range("N32").Value = -30572.334
$ws.Range("H130").Value = 59950
$ws.Range("J130").Value = 59950
$ws.Range("L130").Value = 59950
$ws.Range("N130").Value = -69990

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 8080.636
$ws.Range("I82").Value = 5061.4
$ws.Range("J82").Value = 38273
$ws.Range("K82").Value = 5061.4
$ws.Range("L82").Value = 38273
$ws.Range("M82").Value = -4678.4
$ws.Range("N82").Value = -39039
$ws.Range("H85").Value = 8080.636
$ws.Range("I85").Value = 5061.4
$ws.Range("J85").Value = 38273
$ws.Range("K85").Value = 5061.4
$ws.Range("L85").Value = 38273
$ws.Range("M85").Value = -3735.4
$ws.Range("N85").Value = -40925
$ws.Range("H105").Value = 4926.643
$ws.Range("I105").Value = 4717.8823
$ws.Range("J105").Value = 5813.875
$ws.Range("K105").Value = 4717.8823
$ws.Range("L105").Value = 5813.875
$ws.Range("M105").Value = -2970.8823
$ws.Range("N105").Value = -9307.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3366.4062
$ws.Range("I31").Value = 1297.2858
$ws.Range("J31").Value = 4975.722
$ws.Range("K31").Value = 1297.2858
$ws.Range("L31").Value = 4975.722
$ws.Range("M31").Value = -1002.2858
$ws.Range("N31").Value = -5565.722
$ws.Range("H34").Value = 3366.4062
$ws.Range("I34").Value = 1297.2858
$ws.Range("J34").Value = 4975.722
$ws.Range("K34").Value = 1297.2858
$ws.Range("L34").Value = 4975.722
$ws.Range("M34").Value = -1095.2858
$ws.Range("N34").Value = -5379.722
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 10000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -9572
$ws.Range("N41").ClearContents()
$ws.Range("H51").Value = 25221.777
$ws.Range("I51").Value = 15000
$ws.Range("J51").Value = 26499.5
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 26499.5
$ws.Range("M51").Value = -14264
$ws.Range("N51").Value = -27971.5
$ws.Range("H59").Value = 19727.273
$ws.Range("J59").Value = 19900
$ws.Range("L59").Value = 19900
$ws.Range("N59").Value = -22190
$ws.Range("H60").Value = 19287.875
$ws.Range("I60").Value = 13700
$ws.Range("J60").Value = 28601
$ws.Range("K60").Value = 13700
$ws.Range("L60").Value = 28601
$ws.Range("M60").Value = -13189
$ws.Range("N60").Value = -29623
$ws.Range("H61").Value = 25221.777
$ws.Range("I61").Value = 15000
$ws.Range("J61").Value = 26499.5
$ws.Range("K61").Value = 15000
$ws.Range("L61").Value = 26499.5
$ws.Range("M61").Value = -14652
$ws.Range("N61").Value = -27195.5
$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -17246
$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -56232
$ws.Range("H99").Value = 2850
$ws.Range("I99").Value = 2820
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2820
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1322
$ws.Range("N99").Value = -5996
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 74000
$ws.Range("J125").Value = 74000
$ws.Range("L125").Value = 74000
$ws.Range("N125").Value = -78920
$ws.Range("H126").Value = 2850
$ws.Range("I126").Value = 2820
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8460
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5990
$ws.Range("N126").Value = -13940
$ws.Range("H127").Value = 78933.336
$ws.Range("J127").Value = 78933.336
$ws.Range("L127").Value = 78933.336
$ws.Range("N127").Value = -88853.336

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2648338
$ws.Range("I5").Value = 625.3143
$ws.Range("J5").Value = 5957978.5
$ws.Range("K5").Value = 1875.9429
$ws.Range("L5").Value = 17873935.5
$ws.Range("M5").Value = -1763.9429
$ws.Range("N5").Value = -17874159.5
$ws.Range("H113").Value = 777.2143
$ws.Range("I113").Value = 786.96295
$ws.Range("K113").Value = 2360.88885
$ws.Range("M113").Value = -190.8888499999998
$ws.Range("H122").Value = 1313.8206
$ws.Range("I122").Value = 587.9231
$ws.Range("J122").Value = 1676.7693
$ws.Range("K122").Value = 5291.3079
$ws.Range("L122").Value = 15090.9237
$ws.Range("M122").Value = -2841.3079
$ws.Range("N122").Value = -19990.9237
$ws.Range("H135").Value = 2648338
$ws.Range("I135").Value = 625.3143
$ws.Range("J135").Value = 5957978.5
$ws.Range("K135").Value = 5627.8287
$ws.Range("L135").Value = 53621806.5
$ws.Range("M135").Value = -3092.8287
$ws.Range("N135").Value = -53626876.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 11000
$ws.Range("H70").Value = 5530.978
$ws.Range("I70").Value = 5116.64
$ws.Range("J70").Value = 6024.2383
$ws.Range("K70").Value = 5116.64
$ws.Range("L70").Value = 6024.2383
$ws.Range("M70").Value = -4846.64
$ws.Range("N70").Value = -6564.2383
$ws.Range("H73").Value = 5530.978
$ws.Range("I73").Value = 5116.64
$ws.Range("J73").Value = 6024.2383
$ws.Range("K73").Value = 5116.64
$ws.Range("L73").Value = 6024.2383
$ws.Range("M73").Value = -4180.64
$ws.Range("N73").Value = -7896.2383
$ws.Range("H80").Value = 6043.75
$ws.Range("I80").Value = 9631.25
$ws.Range("J80").Value = 4250
$ws.Range("K80").Value = 9631.25
$ws.Range("L80").Value = 4250
$ws.Range("M80").Value = -8633.25
$ws.Range("N80").Value = -6246
$ws.Range("H83").Value = 6043.75
$ws.Range("I83").Value = 9631.25
$ws.Range("J83").Value = 4250
$ws.Range("K83").Value = 48156.25
$ws.Range("L83").Value = 21250
$ws.Range("M83").Value = -43164.25
$ws.Range("N83").Value = -31234

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4683
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 498.75
$ws.Range("I33").Value = 498.75
$ws.Range("K33").Value = 498.75
$ws.Range("M33").Value = -208.75
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H39").Value = 3068.5715
$ws.Range("J39").Value = 3068.5715
$ws.Range("L39").Value = 3068.5715
$ws.Range("N39").Value = -3988.5715

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 17920
$ws.Range("J58").Value = 17920
$ws.Range("L58").Value = 17920
$ws.Range("N58").Value = -18536
$ws.Range("H109").Value = 63574
$ws.Range("J109").Value = 63574
$ws.Range("L109").Value = 63574
$ws.Range("N109").Value = -66348
